# "Ajout de l'enregistrement vocal"
#
# The questionnaire sheet ("sub-03_task-work_questionnaire_beh") gets its
# real logged incidents filled in: row 2's Timecode is corrected and its
# blank columns are populated, and seven more event rows (3-9) are appended
# below it (run-01, alt-tabs, app freezes, a notification interruption…).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=ID, B=Path, C=Timecode, D=Parameter, E=ID Cible,
#          F=nature_incident, G=responsabilite, H=tache, I=importance,
#          J=description_incident, K=concentration, L=distraction,
#          M=nature_distraction, N=fatigue, O=difficulte
$rows = @(
    @(1, 'run-01', 647,  0, '--', 'Action',       'Utilisateur (Moi)', 'Code : Python, Java, Html ', 'Insignifiante', 'alt tab',              'Neutre', '',    '',             'Neutre',         ''),
    @(2, 'run-01', 687,  3, '--', '',              '',                  '',                           '',              'alt tab',              '',       '',    '',             '',               ''),
    @(3, 'run-01', 790,  0, '--', 'Matériel',      'Système (Machine)', 'Code : Python, Java, Html ', 'Insignifiante', "L'application freeze", 'Neutre', '',    '',             'Neutre',         ''),
    @(4, 'run-01', 1150, 0, '--', 'Matériel',      'Système (Machine)', 'Code : Python, Java, Html ', 'Insignifiante', "l'application freeze", 'Neutre', '',    '',             'Neutre',         ''),
    @(5, 'run-01', 3916, 3, '--', '',              '',                  '',                           '',              'alt tab',              '',       '',    '',             '',               ''),
    @(6, 'run-01', 4020, 3, '--', '',              '',                  '',                           '',              'alt tab',              '',       '',    '',             '',               ''),
    @(7, 'run-01', 4676, 0, '--', 'Perturbation',  'Utilisateur (Moi)', 'Code : Python, Java, Html ', 'Insignifiante', 'Notification',         'Faible', 'Oui', 'notification', 'Plutôt fatigué', 'Moyenne'),
    @(8, 'run-01', 5003, 3, '--', '',              '',                  '',                           '',              'alt tab',              '',       '',    '',             '',               '')
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
